$d = $word.ActiveDocument
$n = $d.Paragraphs.Count

# The last two body paragraphs (both styled "TS"):
#   - second-to-last: empty TS paragraph -> gains rPr/rFonts (cs="Times New Roman")
#   - last: TS paragraph with a numbered list (numId 2) containing a single
#     space run -> loses the list numbering + run text, gains a 720-twip
#     left indent and rPr/rFonts (cs="Times New Roman")
$pLast = $d.Paragraphs.Item($n)
$pPrev = $d.Paragraphs.Item($n - 1)

# Re-applying the paragraph's own style strips the direct w:numPr list
# formatting that was attached to it (do this before removing the run so
# the paragraph mark isn't the only content, which avoids an extraneous
# rsid stamp on the paragraph mark's run properties).
$pLast.Style = "TS"

# Remove the lone space run from the last paragraph, keeping its mark.
$runRange = $d.Range($pLast.Range.Start, $pLast.Range.End - 1)
$runRange.Delete()

# Give it the 720-twip (36pt) left indent the list used to provide visually.
$pLast.LeftIndent = 36

# Stamp both paragraph marks with the Times New Roman complex-script font.
$pLast.Range.Font.NameBi = "Times New Roman"
$pPrev.Range.Font.NameBi = "Times New Roman"
